$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / used range implicitly by writing new rows 14-17

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf1"
$ws.Cells.Item(2,3).Value = "Cd44"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.5247423333333333
$ws.Cells.Item(2,8).Value = 1.574227
$ws.Cells.Item(2,9).Value = 0.044891155074209
$ws.Cells.Item(2,10).Value = 0.044891155074209
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 31.82741333333333
$ws.Cells.Item(2,14).Value = 95.48223999999999
$ws.Cells.Item(2,15).Value = 0.114390792932228
$ws.Cells.Item(2,16).Value = 0.114390792932228
$ws.Cells.Item(2,17).Value = 16.70119113649778
$ws.Cells.Item(2,18).Value = 150.31072022848
$ws.Cells.Item(2,19).Value = 0.005135134824582378
$ws.Cells.Item(2,20).Value = 0.00513513482458238

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf1"
$ws.Cells.Item(3,3).Value = "Cd44"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.5247423333333333
$ws.Cells.Item(3,8).Value = 1.574227
$ws.Cells.Item(3,9).Value = 0.044891155074209
$ws.Cells.Item(3,10).Value = 0.044891155074209
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 85.46317833333335
$ws.Cells.Item(3,14).Value = 256.389535
$ws.Cells.Item(3,15).Value = 0.307162904935779
$ws.Cells.Item(3,16).Value = 0.307162904935779
$ws.Cells.Item(3,17).Value = 44.84614761271612
$ws.Cells.Item(3,18).Value = 403.615328514445
$ws.Cells.Item(3,19).Value = 0.01378889759851657
$ws.Cells.Item(3,20).Value = 0.01378889759851657

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf1"
$ws.Cells.Item(4,3).Value = "Cd44"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.5247423333333333
$ws.Cells.Item(4,8).Value = 1.574227
$ws.Cells.Item(4,9).Value = 0.044891155074209
$ws.Cells.Item(4,10).Value = 0.044891155074209
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 122.2478306666667
$ws.Cells.Item(4,14).Value = 366.743492
$ws.Cells.Item(4,15).Value = 0.4393704929064738
$ws.Cells.Item(4,16).Value = 0.4393704929064738
$ws.Cells.Item(4,17).Value = 64.1486119089649
$ws.Cells.Item(4,18).Value = 577.337507180684
$ws.Cells.Item(4,19).Value = 0.01972384893209616
$ws.Cells.Item(4,20).Value = 0.01972384893209616

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Fgf1"
$ws.Cells.Item(5,3).Value = "Cd44"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.5247423333333333
$ws.Cells.Item(5,8).Value = 1.574227
$ws.Cells.Item(5,9).Value = 0.044891155074209
$ws.Cells.Item(5,10).Value = 0.044891155074209
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 38.69562533333333
$ws.Cells.Item(5,14).Value = 116.086876
$ws.Cells.Item(5,15).Value = 0.1390758092255191
$ws.Cells.Item(5,16).Value = 0.1390758092255191
$ws.Cells.Item(5,17).Value = 20.30523272720578
$ws.Cells.Item(5,18).Value = 182.747094544852
$ws.Cells.Item(5,19).Value = 0.006243273719013886
$ws.Cells.Item(5,20).Value = 0.006243273719013886

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf1"
$ws.Cells.Item(6,3).Value = "Cd44"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 4.531154666666667
$ws.Cells.Item(6,8).Value = 13.593464
$ws.Cells.Item(6,9).Value = 0.3876355191593572
$ws.Cells.Item(6,10).Value = 0.3876355191593572
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 31.82741333333333
$ws.Cells.Item(6,14).Value = 95.48223999999999
$ws.Cells.Item(6,15).Value = 0.114390792932228
$ws.Cells.Item(6,16).Value = 0.114390792932228
$ws.Cells.Item(6,17).Value = 144.2149324532622
$ws.Cells.Item(6,18).Value = 1297.93439207936
$ws.Cells.Item(6,19).Value = 0.04434193440533474
$ws.Cells.Item(6,20).Value = 0.04434193440533474

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf1"
$ws.Cells.Item(7,3).Value = "Cd44"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 4.531154666666667
$ws.Cells.Item(7,8).Value = 13.593464
$ws.Cells.Item(7,9).Value = 0.3876355191593572
$ws.Cells.Item(7,10).Value = 0.3876355191593572
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 85.46317833333335
$ws.Cells.Item(7,14).Value = 256.389535
$ws.Cells.Item(7,15).Value = 0.307162904935779
$ws.Cells.Item(7,16).Value = 0.307162904935779
$ws.Cells.Item(7,17).Value = 387.246879333249
$ws.Cells.Item(7,18).Value = 3485.22191399924
$ws.Cells.Item(7,19).Value = 0.119067252121277
$ws.Cells.Item(7,20).Value = 0.119067252121277

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Fgf1"
$ws.Cells.Item(8,3).Value = "Cd44"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.531154666666667
$ws.Cells.Item(8,8).Value = 13.593464
$ws.Cells.Item(8,9).Value = 0.3876355191593572
$ws.Cells.Item(8,10).Value = 0.3876355191593572
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 122.2478306666667
$ws.Cells.Item(8,14).Value = 366.743492
$ws.Cells.Item(8,15).Value = 0.4393704929064738
$ws.Cells.Item(8,16).Value = 0.4393704929064738
$ws.Cells.Item(8,17).Value = 553.9238284151431
$ws.Cells.Item(8,18).Value = 4985.314455736288
$ws.Cells.Item(8,19).Value = 0.1703156091211037
$ws.Cells.Item(8,20).Value = 0.1703156091211037

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Fgf1"
$ws.Cells.Item(9,3).Value = "Cd44"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.531154666666667
$ws.Cells.Item(9,8).Value = 13.593464
$ws.Cells.Item(9,9).Value = 0.3876355191593572
$ws.Cells.Item(9,10).Value = 0.3876355191593572
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 38.69562533333333
$ws.Cells.Item(9,14).Value = 116.086876
$ws.Cells.Item(9,15).Value = 0.1390758092255191
$ws.Cells.Item(9,16).Value = 0.1390758092255191
$ws.Cells.Item(9,17).Value = 175.3358633087182
$ws.Cells.Item(9,18).Value = 1578.022769778464
$ws.Cells.Item(9,19).Value = 0.05391072351164183
$ws.Cells.Item(9,20).Value = 0.05391072351164182

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Fgf1"
$ws.Cells.Item(10,3).Value = "Cd44"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.08000233333333333
$ws.Cells.Item(10,8).Value = 0.240007
$ws.Cells.Item(10,9).Value = 0.006844115528380393
$ws.Cells.Item(10,10).Value = 0.006844115528380393
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 31.82741333333333
$ws.Cells.Item(10,14).Value = 95.48223999999999
$ws.Cells.Item(10,15).Value = 0.114390792932228
$ws.Cells.Item(10,16).Value = 0.114390792932228
$ws.Cells.Item(10,17).Value = 2.546267330631111
$ws.Cells.Item(10,18).Value = 22.91640597568
$ws.Cells.Item(10,19).Value = 0.0007829038022112078
$ws.Cells.Item(10,20).Value = 0.000782903802211208

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Fgf1"
$ws.Cells.Item(11,3).Value = "Cd44"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.08000233333333333
$ws.Cells.Item(11,8).Value = 0.240007
$ws.Cells.Item(11,9).Value = 0.006844115528380393
$ws.Cells.Item(11,10).Value = 0.006844115528380393
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 85.46317833333335
$ws.Cells.Item(11,14).Value = 256.389535
$ws.Cells.Item(11,15).Value = 0.307162904935779
$ws.Cells.Item(11,16).Value = 0.307162904935779
$ws.Cells.Item(11,17).Value = 6.837253680749445
$ws.Cells.Item(11,18).Value = 61.53528312674501
$ws.Cells.Item(11,19).Value = 0.002102258407413396
$ws.Cells.Item(11,20).Value = 0.002102258407413396

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Fgf1"
$ws.Cells.Item(12,3).Value = "Cd44"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.08000233333333333
$ws.Cells.Item(12,8).Value = 0.240007
$ws.Cells.Item(12,9).Value = 0.006844115528380393
$ws.Cells.Item(12,10).Value = 0.006844115528380393
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 122.2478306666667
$ws.Cells.Item(12,14).Value = 366.743492
$ws.Cells.Item(12,15).Value = 0.4393704929064738
$ws.Cells.Item(12,16).Value = 0.4393704929064738
$ws.Cells.Item(12,17).Value = 9.780111698271556
$ws.Cells.Item(12,18).Value = 88.021005284444
$ws.Cells.Item(12,19).Value = 0.003007102413213344
$ws.Cells.Item(12,20).Value = 0.003007102413213345

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Fgf1"
$ws.Cells.Item(13,3).Value = "Cd44"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.08000233333333333
$ws.Cells.Item(13,8).Value = 0.240007
$ws.Cells.Item(13,9).Value = 0.006844115528380393
$ws.Cells.Item(13,10).Value = 0.006844115528380393
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 38.69562533333333
$ws.Cells.Item(13,14).Value = 116.086876
$ws.Cells.Item(13,15).Value = 0.1390758092255191
$ws.Cells.Item(13,16).Value = 0.1390758092255191
$ws.Cells.Item(13,17).Value = 3.095740316459111
$ws.Cells.Item(13,18).Value = 27.861662848132
$ws.Cells.Item(13,19).Value = 0.0009518509055424446
$ws.Cells.Item(13,20).Value = 0.0009518509055424446

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Fgf1"
$ws.Cells.Item(14,3).Value = "Cd44"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 6.553315
$ws.Cells.Item(14,8).Value = 19.659945
$ws.Cells.Item(14,9).Value = 0.5606292102380533
$ws.Cells.Item(14,10).Value = 0.5606292102380533
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 31.82741333333333
$ws.Cells.Item(14,14).Value = 95.48223999999999
$ws.Cells.Item(14,15).Value = 0.114390792932228
$ws.Cells.Item(14,16).Value = 0.114390792932228
$ws.Cells.Item(14,17).Value = 208.5750652085333
$ws.Cells.Item(14,18).Value = 1877.1755868768
$ws.Cells.Item(14,19).Value = 0.06413081990009968
$ws.Cells.Item(14,20).Value = 0.06413081990009969

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Fgf1"
$ws.Cells.Item(15,3).Value = "Cd44"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 6.553315
$ws.Cells.Item(15,8).Value = 19.659945
$ws.Cells.Item(15,9).Value = 0.5606292102380533
$ws.Cells.Item(15,10).Value = 0.5606292102380533
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 85.46317833333335
$ws.Cells.Item(15,14).Value = 256.389535
$ws.Cells.Item(15,15).Value = 0.307162904935779
$ws.Cells.Item(15,16).Value = 0.307162904935779
$ws.Cells.Item(15,17).Value = 560.0671285195084
$ws.Cells.Item(15,18).Value = 5040.604156675576
$ws.Cells.Item(15,19).Value = 0.1722044968085721
$ws.Cells.Item(15,20).Value = 0.1722044968085721

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Fgf1"
$ws.Cells.Item(16,3).Value = "Cd44"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 6.553315
$ws.Cells.Item(16,8).Value = 19.659945
$ws.Cells.Item(16,9).Value = 0.5606292102380533
$ws.Cells.Item(16,10).Value = 0.5606292102380533
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 122.2478306666667
$ws.Cells.Item(16,14).Value = 366.743492
$ws.Cells.Item(16,15).Value = 0.4393704929064738
$ws.Cells.Item(16,16).Value = 0.4393704929064738
$ws.Cells.Item(16,17).Value = 801.1285424253267
$ws.Cells.Item(16,18).Value = 7210.15688182794
$ws.Cells.Item(16,19).Value = 0.2463239324400606
$ws.Cells.Item(16,20).Value = 0.2463239324400607

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Fgf1"
$ws.Cells.Item(17,3).Value = "Cd44"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 6.553315
$ws.Cells.Item(17,8).Value = 19.659945
$ws.Cells.Item(17,9).Value = 0.5606292102380533
$ws.Cells.Item(17,10).Value = 0.5606292102380533
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 38.69562533333333
$ws.Cells.Item(17,14).Value = 116.086876
$ws.Cells.Item(17,15).Value = 0.1390758092255191
$ws.Cells.Item(17,16).Value = 0.1390758092255191
$ws.Cells.Item(17,17).Value = 253.5846219313133
$ws.Cells.Item(17,18).Value = 2282.26159738182
$ws.Cells.Item(17,19).Value = 0.07796996108932096
$ws.Cells.Item(17,20).Value = 0.07796996108932096
